# Refresh the "starttime" placeholder timestamp in column J.
# The MySQL query that produced this export re-ran and every row that
# still carried the previous placeholder start time needs to be bumped
# to the new value returned by the (now-optimized) inquire.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldStartTime = 1586659823.314746
$newStartTime = 1586660079.314404

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("J" + $r)
    $val = $cell.Value2
    if ($val -ne $null -and [Math]::Abs([double]$val - $oldStartTime) -lt 0.0000005) {
        $cell.Value2 = $newStartTime
    }
}
